# daily auto push: 2026-01-30 02:50 UTC
# Insert a new scrape row for 2026/01/30 (Friday) at row 722, pushing all
# subsequent rows (old 722..763) down by one (to 723..764).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a blank row above the current row 722 (which contained the first
# 2026/12/29 entry); that row and everything below it shifts down by one.
$ws.Rows("722:722").Insert()

# Column A holds a literal date-like string (e.g. "2026/01/30"), not a real
# Excel date. Force text formatting before assignment so Excel/COM does not
# auto-convert it into a date serial number, then clear the number format
# back to the sheet default so no stray style is left behind.
$ws.Range("A722").NumberFormat = "@"
$ws.Range("A722").Value = "2026/01/30"
$ws.Range("A722").Style = "Normal"

$ws.Range("B722").Value = "金"
$ws.Range("C722").Value = 9
$ws.Range("D722").Value = 22
